$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I33").Value = 1135
$ws.Range("H33").Value = 1135
$ws.Range("K33").Value = 1135
$ws.Range("M33").Value = -906
$ws.Range("H45").Value = 4969
$ws.Range("L45").Value = 14907
$ws.Range("J45").Value = 4969
$ws.Range("N45").Value = -15291
$ws.Range("I69").Value = 0
$ws.Range("H69").Value = 10000
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("L69").Value = 30000
$ws.Range("J69").Value = 10000
$ws.Range("N69").Value = -31748
$ws.Range("I72").Value = 0
$ws.Range("H72").Value = 10000
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("L72").Value = 90000
$ws.Range("J72").Value = 10000
$ws.Range("N72").Value = -98736
$ws.Range("I74").Value = 10577.857
$ws.Range("H74").Value = 10622.941
$ws.Range("K74").Value = 10577.857
$ws.Range("M74").Value = -9641.857
$ws.Range("L74").Value = 10833.333
$ws.Range("J74").Value = 10833.333
$ws.Range("N74").Value = -12705.333
$ws.Range("I76").Value = 6226.4443
$ws.Range("K76").Value = 6226.4443
$ws.Range("M76").Value = -5911.4443
$ws.Range("L76").Value = 125006620
$ws.Range("J76").Value = 125006620
$ws.Range("N76").Value = -125007250
$ws.Range("I77").Value = 10577.857
$ws.Range("H77").Value = 10622.941
$ws.Range("K77").Value = 52889.285
$ws.Range("M77").Value = -48209.285
$ws.Range("L77").Value = 54166.665
$ws.Range("J77").Value = 10833.333
$ws.Range("N77").Value = -63526.665
$ws.Range("I79").Value = 6226.4443
$ws.Range("K79").Value = 6226.4443
$ws.Range("M79").Value = -5134.4443
$ws.Range("L79").Value = 125006620
$ws.Range("J79").Value = 125006620
$ws.Range("N79").Value = -125008804
$ws.Range("I80").Value = 8446.385
$ws.Range("H80").Value = 5350.909
$ws.Range("K80").Value = 25339.155
$ws.Range("M80").Value = -24341.155
$ws.Range("I82").Value = 1995
$ws.Range("H82").Value = 1995
$ws.Range("K82").Value = 5985
$ws.Range("M82").Value = -5579
$ws.Range("I83").Value = 8446.385
$ws.Range("H83").Value = 5350.909
$ws.Range("K83").Value = 76017.465
$ws.Range("M83").Value = -71025.465
$ws.Range("I85").Value = 1995
$ws.Range("H85").Value = 1995
$ws.Range("K85").Value = 5985
$ws.Range("M85").Value = -4581
$ws.Range("I86").Value = 2120.625
$ws.Range("H86").Value = 2773370.2
$ws.Range("K86").Value = 2120.625
$ws.Range("M86").Value = -997.625
$ws.Range("L86").Value = 4788824.5
$ws.Range("J86").Value = 4788824.5
$ws.Range("N86").Value = -4791070.5
$ws.Range("I88").Value = 699
$ws.Range("H88").Value = 1841.8334
$ws.Range("K88").Value = 699
$ws.Range("M88").Value = -293
$ws.Range("L88").Value = 2070.4
$ws.Range("J88").Value = 2070.4
$ws.Range("N88").Value = -2882.4
$ws.Range("I89").Value = 2120.625
$ws.Range("H89").Value = 2773370.2
$ws.Range("K89").Value = 10603.125
$ws.Range("M89").Value = -4987.125
$ws.Range("L89").Value = 23944122.5
$ws.Range("J89").Value = 4788824.5
$ws.Range("N89").Value = -23955354.5
$ws.Range("I91").Value = 699
$ws.Range("H91").Value = 1841.8334
$ws.Range("K91").Value = 699
$ws.Range("M91").Value = 705
$ws.Range("L91").Value = 2070.4
$ws.Range("J91").Value = 2070.4
$ws.Range("N91").Value = -4878.4
$ws.Range("I92").Value = 449
$ws.Range("H92").Value = 449
$ws.Range("K92").Value = 449
$ws.Range("M92").Value = 799
$ws.Range("I94").Value = 1620.7778
$ws.Range("H94").Value = 2069.4167
$ws.Range("K94").Value = 1620.7778
$ws.Range("M94").Value = -1169.7778
$ws.Range("L94").Value = 3415.3333
$ws.Range("J94").Value = 3415.3333
$ws.Range("N94").Value = -4317.3333
$ws.Range("I96").Value = 482
$ws.Range("H96").Value = 533.8
$ws.Range("K96").Value = 1446
$ws.Range("M96").Value = -73
$ws.Range("L96").Value = 3000
$ws.Range("J96").Value = 1000
$ws.Range("N96").Value = -5746
$ws.Range("I98").Value = 870.3889
$ws.Range("H98").Value = 1087.7368
$ws.Range("K98").Value = 870.3889
$ws.Range("M98").Value = 627.6111
$ws.Range("I100").Value = 1743.7778
$ws.Range("H100").Value = 8784.5
$ws.Range("K100").Value = 1743.7778
$ws.Range("M100").Value = -1202.7778
$ws.Range("L100").Value = 14545.091
$ws.Range("J100").Value = 14545.091
$ws.Range("N100").Value = -15627.091
$ws.Range("I101").Value = 675.1667
$ws.Range("H101").Value = 721.5714
$ws.Range("K101").Value = 2025.5001
$ws.Range("M101").Value = -403.5001
$ws.Range("L101").Value = 3000
$ws.Range("J101").Value = 1000
$ws.Range("N101").Value = -6244
$ws.Range("I104").Value = 1013.3333
$ws.Range("H104").Value = 1013.3333
$ws.Range("K104").Value = 3039.9999
$ws.Range("M104").Value = -1292.9999
$ws.Range("I107").Value = 42639.5
$ws.Range("H107").Value = 39416.348
$ws.Range("K107").Value = 42639.5
$ws.Range("M107").Value = -40719.5
$ws.Range("I113").Value = 0
$ws.Range("H113").Value = 3000
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("L113").Value = 3000
$ws.Range("J113").Value = 3000
$ws.Range("N113").Value = -9508
$ws.Range("I122").Value = 870.3889
$ws.Range("H122").Value = 1087.7368
$ws.Range("K122").Value = 2611.1667
$ws.Range("M122").Value = -161.1667000000002
$ws.Range("H128").Value = 89556.84
$ws.Range("L128").Value = 89556.84
$ws.Range("J128").Value = 89556.84
$ws.Range("N128").Value = -99516.84
$ws.Range("H130").Value = 74786
$ws.Range("L130").Value = 74786
$ws.Range("J130").Value = 74786
$ws.Range("N130").Value = -84826
$ws.Range("I132").Value = 1728.0444
$ws.Range("H132").Value = 9632.588
$ws.Range("K132").Value = 5184.1332
$ws.Range("M132").Value = -2654.1332
$ws.Range("I137").Value = 5869.227
$ws.Range("H137").Value = 5800.2593
$ws.Range("K137").Value = 17607.681
$ws.Range("M137").Value = -15057.681
$ws.Range("L137").Value = 16490.4
$ws.Range("J137").Value = 5496.8
$ws.Range("N137").Value = -21590.4
$ws.Range("I138").Value = 3101.1875
$ws.Range("H138").Value = 4942.153
$ws.Range("K138").Value = 9303.5625
$ws.Range("M138").Value = -4163.5625
$ws.Range("L138").Value = 16404.429
$ws.Range("J138").Value = 5468.143
$ws.Range("N138").Value = -26684.429
$ws.Range("I141").Value = 3788.3
$ws.Range("H141").Value = 4362.5454
$ws.Range("K141").Value = 11364.9
$ws.Range("M141").Value = -6184.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 3526.1267
$ws.Range("H32").Value = 4010.76
$ws.Range("K32").Value = 3526.1267
$ws.Range("M32").Value = -3239.1267
$ws.Range("L32").Value = 12613
$ws.Range("J32").Value = 12613
$ws.Range("N32").Value = -13187
$ws.Range("I61").Value = 2946.8147
$ws.Range("H61").Value = 3191.862
$ws.Range("K61").Value = 2946.8147
$ws.Range("M61").Value = -2734.8147
$ws.Range("I74").Value = 1914.8572
$ws.Range("H74").Value = 1914.8572
$ws.Range("K74").Value = 1914.8572
$ws.Range("M74").Value = -1040.8572
$ws.Range("L74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("I77").Value = 1914.8572
$ws.Range("H77").Value = 1914.8572
$ws.Range("K77").Value = 9574.286
$ws.Range("M77").Value = -5206.286
$ws.Range("L77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("I88").Value = 3760.75
$ws.Range("H88").Value = 3172.6365
$ws.Range("K88").Value = 3760.75
$ws.Range("M88").Value = -3354.75
$ws.Range("L88").Value = 1604.3334
$ws.Range("J88").Value = 1604.3334
$ws.Range("N88").Value = -2416.3334
$ws.Range("I91").Value = 3760.75
$ws.Range("H91").Value = 3172.6365
$ws.Range("K91").Value = 3760.75
$ws.Range("M91").Value = -2356.75
$ws.Range("L91").Value = 1604.3334
$ws.Range("J91").Value = 1604.3334
$ws.Range("N91").Value = -4412.3334
$ws.Range("I97").Value = 1071
$ws.Range("H97").Value = 1022.6111
$ws.Range("K97").Value = 1071
$ws.Range("M97").Value = -575
$ws.Range("L97").Value = 200
$ws.Range("J97").Value = 200
$ws.Range("N97").Value = -1192
$ws.Range("H98").Value = 9000
$ws.Range("L98").Value = 9000
$ws.Range("J98").Value = 9000
$ws.Range("N98").Value = -14990
$ws.Range("I102").Value = 1818
$ws.Range("H102").Value = 2266.5334
$ws.Range("K102").Value = 1818
$ws.Range("M102").Value = -196
$ws.Range("I110").Value = 2819.3
$ws.Range("H110").Value = 2682.6667
$ws.Range("K110").Value = 2819.3
$ws.Range("M110").Value = -774.3000000000002
$ws.Range("I122").Value = 4188.5
$ws.Range("H122").Value = 4869.108
$ws.Range("K122").Value = 12565.5
$ws.Range("M122").Value = -10115.5
$ws.Range("L122").Value = 16541.6835
$ws.Range("J122").Value = 5513.8945
$ws.Range("N122").Value = -21441.6835
$ws.Range("H133").Value = 59999.668
$ws.Range("L133").Value = 59999.668
$ws.Range("J133").Value = 59999.668
$ws.Range("N133").Value = -65059.668
$ws.Range("H134").Value = 52052.25
$ws.Range("L134").Value = 53283.6
$ws.Range("J134").Value = 53283.6
$ws.Range("N134").Value = -63423.6
$ws.Range("I136").Value = 2946.8147
$ws.Range("H136").Value = 3191.862
$ws.Range("K136").Value = 8840.444100000001
$ws.Range("M136").Value = -6290.444100000001
$ws.Range("H137").Value = 73091.8
$ws.Range("L137").Value = 73091.8
$ws.Range("J137").Value = 73091.8
$ws.Range("N137").Value = -83291.8
$ws.Range("H139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("N139").Value = -60280
$ws.Range("H140").Value = 58333.332
$ws.Range("L140").Value = 58333.332
$ws.Range("J140").Value = 58333.332
$ws.Range("N140").Value = -68693.33199999999
$ws.Range("H141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 2430.85
$ws.Range("H20").Value = 2482.9714
$ws.Range("K20").Value = 2430.85
$ws.Range("M20").Value = -2183.85
$ws.Range("L20").Value = 2552.4666
$ws.Range("J20").Value = 2552.4666
$ws.Range("N20").Value = -3046.4666
$ws.Range("I64").Value = 853
$ws.Range("H64").Value = 844.3077
$ws.Range("K64").Value = 853
$ws.Range("M64").Value = -628
$ws.Range("L64").Value = 842.7273
$ws.Range("J64").Value = 842.7273
$ws.Range("N64").Value = -1292.7273
$ws.Range("I67").Value = 853
$ws.Range("H67").Value = 844.3077
$ws.Range("K67").Value = 853
$ws.Range("M67").Value = -73
$ws.Range("L67").Value = 842.7273
$ws.Range("J67").Value = 842.7273
$ws.Range("N67").Value = -2402.7273
$ws.Range("I86").Value = 946796.9
$ws.Range("H86").Value = 656644.6
$ws.Range("K86").Value = 946796.9
$ws.Range("M86").Value = -945673.9
$ws.Range("I89").Value = 946796.9
$ws.Range("H89").Value = 656644.6
$ws.Range("K89").Value = 4733984.5
$ws.Range("M89").Value = -4728368.5
$ws.Range("I94").Value = 2534.0833
$ws.Range("H94").Value = 2187.2666
$ws.Range("K94").Value = 2534.0833
$ws.Range("M94").Value = -2083.0833
$ws.Range("L94").Value = 800
$ws.Range("J94").Value = 800
$ws.Range("N94").Value = -1702
$ws.Range("I99").Value = 2098.8
$ws.Range("H99").Value = 2344.1538
$ws.Range("K99").Value = 2098.8
$ws.Range("M99").Value = -600.8000000000002
$ws.Range("I105").Value = 127682.5
$ws.Range("H105").Value = 113628.78
$ws.Range("K105").Value = 127682.5
$ws.Range("M105").Value = -125935.5
$ws.Range("L105").Value = 1199
$ws.Range("J105").Value = 1199
$ws.Range("N105").Value = -4693
$ws.Range("I107").Value = 2691.2
$ws.Range("H107").Value = 1253433.1
$ws.Range("K107").Value = 2691.2
$ws.Range("M107").Value = -771.1999999999998
$ws.Range("L107").Value = 3338003
$ws.Range("J107").Value = 3338003
$ws.Range("N107").Value = -3341843
$ws.Range("I134").Value = 5759.4814
$ws.Range("H134").Value = 36878.25
$ws.Range("K134").Value = 17278.4442
$ws.Range("M134").Value = -14743.4442
$ws.Range("L134").Value = 614758.8
$ws.Range("J134").Value = 204919.6
$ws.Range("N134").Value = -619828.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 1629
$ws.Range("H31").Value = 52053.332
$ws.Range("K31").Value = 1629
$ws.Range("M31").Value = -1334
$ws.Range("L31").Value = 67810.94
$ws.Range("J31").Value = 67810.94
$ws.Range("N31").Value = -68400.94
$ws.Range("I34").Value = 1629
$ws.Range("H34").Value = 52053.332
$ws.Range("K34").Value = 1629
$ws.Range("M34").Value = -1427
$ws.Range("L34").Value = 67810.94
$ws.Range("J34").Value = 67810.94
$ws.Range("N34").Value = -68214.94
$ws.Range("I58").Value = 5525.5
$ws.Range("H58").Value = 4620.4
$ws.Range("K58").Value = 5525.5
$ws.Range("M58").Value = -5322.5
$ws.Range("L58").Value = 1000
$ws.Range("J58").Value = 1000
$ws.Range("N58").Value = -1406
$ws.Range("I86").Value = 10000.5
$ws.Range("H86").Value = 10277.75
$ws.Range("K86").Value = 10000.5
$ws.Range("M86").Value = -8877.5
$ws.Range("L86").Value = 10555
$ws.Range("J86").Value = 10555
$ws.Range("N86").Value = -12801
$ws.Range("I89").Value = 10000.5
$ws.Range("H89").Value = 10277.75
$ws.Range("K89").Value = 50002.5
$ws.Range("M89").Value = -44386.5
$ws.Range("L89").Value = 52775
$ws.Range("J89").Value = 10555
$ws.Range("N89").Value = -64007
$ws.Range("I105").Value = 743
$ws.Range("H105").Value = 748
$ws.Range("K105").Value = 743
$ws.Range("M105").Value = 1004
$ws.Range("I107").Value = 443.05264
$ws.Range("H107").Value = 499.56
$ws.Range("K107").Value = 443.05264
$ws.Range("M107").Value = 1476.94736
$ws.Range("L107").Value = 678.5
$ws.Range("J107").Value = 678.5
$ws.Range("N107").Value = -4518.5
$ws.Range("I122").Value = 3036.25
$ws.Range("H122").Value = 3857.5417
$ws.Range("K122").Value = 9108.75
$ws.Range("M122").Value = -6658.75
$ws.Range("L122").Value = 16500.375
$ws.Range("J122").Value = 5500.125
$ws.Range("N122").Value = -21400.375
$ws.Range("I132").Value = 1544.6296
$ws.Range("H132").Value = 1816.0588
$ws.Range("K132").Value = 4633.8888
$ws.Range("M132").Value = -2103.8888
$ws.Range("I134").Value = 2544.1562
$ws.Range("H134").Value = 305497.38
$ws.Range("K134").Value = 7632.4686
$ws.Range("M134").Value = -5097.4686
$ws.Range("I136").Value = 5525.5
$ws.Range("H136").Value = 4620.4
$ws.Range("K136").Value = 16576.5
$ws.Range("M136").Value = -14026.5
$ws.Range("L136").Value = 3000
$ws.Range("J136").Value = 1000
$ws.Range("N136").Value = -8100
$ws.Range("H141").Value = 309436.5
$ws.Range("L141").Value = 369656.62
$ws.Range("J141").Value = 369656.62
$ws.Range("N141").Value = -380016.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 101017
$ws.Range("H5").Value = 1200904.5
$ws.Range("K5").Value = 303051
$ws.Range("M5").Value = -302939
$ws.Range("I11").Value = 1666966.8
$ws.Range("H11").Value = 952960
$ws.Range("K11").Value = 5000900.4
$ws.Range("M11").Value = -5000760.4
$ws.Range("I29").Value = 201
$ws.Range("H29").Value = 201
$ws.Range("K29").Value = 603
$ws.Range("M29").Value = -326
$ws.Range("I31").Value = 1000
$ws.Range("H31").Value = 1000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2712
$ws.Range("I56").Value = 6499.75
$ws.Range("H56").Value = 6499.75
$ws.Range("K56").Value = 6499.75
$ws.Range("M56").Value = -5969.75
$ws.Range("I122").Value = 775.6
$ws.Range("H122").Value = 30629.383
$ws.Range("K122").Value = 6980.400000000001
$ws.Range("M122").Value = -4530.400000000001
$ws.Range("L122").Value = 487783.404
$ws.Range("J122").Value = 54198.156
$ws.Range("N122").Value = -492683.404
$ws.Range("I129").Value = 456
$ws.Range("H129").Value = 67913.87
$ws.Range("K129").Value = 1368
$ws.Range("M129").Value = 3632
$ws.Range("L129").Value = 435025.74
$ws.Range("J129").Value = 145008.58
$ws.Range("N129").Value = -445025.74
$ws.Range("I132").Value = 100984.5
$ws.Range("H132").Value = 444313.53
$ws.Range("K132").Value = 908860.5
$ws.Range("M132").Value = -906330.5
$ws.Range("L132").Value = 6058796.040000001
$ws.Range("J132").Value = 673199.5600000001
$ws.Range("N132").Value = -6063856.040000001
$ws.Range("I135").Value = 101017
$ws.Range("H135").Value = 1200904.5
$ws.Range("K135").Value = 909153
$ws.Range("M135").Value = -906618
$ws.Range("I137").Value = 3268.889
$ws.Range("H137").Value = 3142
$ws.Range("K137").Value = 9806.667000000001
$ws.Range("M137").Value = -4706.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I59").Value = 24109
$ws.Range("H59").Value = 24109.75
$ws.Range("K59").Value = 24109
$ws.Range("M59").Value = -23526
$ws.Range("I70").Value = 8256.143
$ws.Range("H70").Value = 10191.77
$ws.Range("K70").Value = 8256.143
$ws.Range("M70").Value = -7986.143
$ws.Range("I73").Value = 8256.143
$ws.Range("H73").Value = 10191.77
$ws.Range("K73").Value = 8256.143
$ws.Range("M73").Value = -7320.143
$ws.Range("I80").Value = 913459.75
$ws.Range("H80").Value = 1433825.1
$ws.Range("K80").Value = 913459.75
$ws.Range("M80").Value = -912461.75
$ws.Range("L80").Value = 3341831.8
$ws.Range("J80").Value = 3341831.8
$ws.Range("N80").Value = -3343827.8
$ws.Range("I83").Value = 913459.75
$ws.Range("H83").Value = 1433825.1
$ws.Range("K83").Value = 4567298.75
$ws.Range("M83").Value = -4562306.75
$ws.Range("L83").Value = 16709159
$ws.Range("J83").Value = 3341831.8
$ws.Range("N83").Value = -16719143
$ws.Range("I122").Value = 2153.7273
$ws.Range("H122").Value = 3631.5264
$ws.Range("K122").Value = 6461.1819
$ws.Range("M122").Value = -4011.1819
$ws.Range("L122").Value = 16990.5
$ws.Range("J122").Value = 5663.5
$ws.Range("N122").Value = -21890.5
$ws.Range("I132").Value = 6259.3706
$ws.Range("H132").Value = 38483.645
$ws.Range("K132").Value = 18778.1118
$ws.Range("M132").Value = -16248.1118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 6603.4287
$ws.Range("H7").Value = 6991.9
$ws.Range("K7").Value = 6603.4287
$ws.Range("M7").Value = -6491.4287
$ws.Range("I61").Value = 10761.4
$ws.Range("H61").Value = 8461.200000000001
$ws.Range("K61").Value = 10761.4
$ws.Range("M61").Value = -10559.4
$ws.Range("L61").Value = 6161
$ws.Range("J61").Value = 6161
$ws.Range("N61").Value = -6565
$ws.Range("H68").Value = 201647.4
$ws.Range("L68").Value = 251698.25
$ws.Range("J68").Value = 251698.25
$ws.Range("N68").Value = -253196.25
$ws.Range("H71").Value = 201647.4
$ws.Range("L71").Value = 1258491.25
$ws.Range("J71").Value = 251698.25
$ws.Range("N71").Value = -1265979.25
$ws.Range("I82").Value = 1911.5
$ws.Range("H82").Value = 1707.1111
$ws.Range("K82").Value = 1911.5
$ws.Range("M82").Value = -1550.5
$ws.Range("L82").Value = 1298.3334
$ws.Range("J82").Value = 1298.3334
$ws.Range("N82").Value = -2020.3334
$ws.Range("I85").Value = 1911.5
$ws.Range("H85").Value = 1707.1111
$ws.Range("K85").Value = 1911.5
$ws.Range("M85").Value = -663.5
$ws.Range("L85").Value = 1298.3334
$ws.Range("J85").Value = 1298.3334
$ws.Range("N85").Value = -3794.3334
$ws.Range("I93").Value = 100001560
$ws.Range("H93").Value = 62502532
$ws.Range("K93").Value = 100001560
$ws.Range("M93").Value = -100000312
$ws.Range("I113").Value = 10761.4
$ws.Range("H113").Value = 8461.200000000001
$ws.Range("K113").Value = 10761.4
$ws.Range("M113").Value = -8591.4
$ws.Range("L113").Value = 6161
$ws.Range("J113").Value = 6161
$ws.Range("N113").Value = -10501
$ws.Range("H122").Value = 6386.75
$ws.Range("L122").Value = 20337
$ws.Range("J122").Value = 6779
$ws.Range("N122").Value = -25237
$ws.Range("I126").Value = 6603.4287
$ws.Range("H126").Value = 6991.9
$ws.Range("K126").Value = 19810.2861
$ws.Range("M126").Value = -17340.2861
$ws.Range("H130").Value = 90214.5
$ws.Range("L130").Value = 90214.5
$ws.Range("J130").Value = 90214.5
$ws.Range("N130").Value = -100254.5
$ws.Range("I132").Value = 5271.364
$ws.Range("H132").Value = 6351.911
$ws.Range("K132").Value = 15814.092
$ws.Range("M132").Value = -13284.092
$ws.Range("H135").Value = 67949.086
$ws.Range("L135").Value = 67949.086
$ws.Range("J135").Value = 67949.086
$ws.Range("N135").Value = -78089.086
$ws.Range("I136").Value = 480456.44
$ws.Range("H136").Value = 292208.44
$ws.Range("K136").Value = 1441369.32
$ws.Range("M136").Value = -1438819.32
$ws.Range("L136").Value = 29509.287
$ws.Range("J136").Value = 9836.429
$ws.Range("N136").Value = -34609.287
$ws.Range("H139").Value = 52458.332
$ws.Range("L139").Value = 52458.332
$ws.Range("J139").Value = 52458.332
$ws.Range("N139").Value = -62738.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 207758.6
$ws.Range("H62").Value = 98798.45
$ws.Range("K62").Value = 207758.6
$ws.Range("M62").Value = -207134.6
$ws.Range("I65").Value = 207758.6
$ws.Range("H65").Value = 98798.45
$ws.Range("K65").Value = 1038793
$ws.Range("M65").Value = -1035673
$ws.Range("H74").Value = 9213
$ws.Range("L74").Value = 9644.5
$ws.Range("J74").Value = 9644.5
$ws.Range("N74").Value = -11516.5
$ws.Range("H77").Value = 9213
$ws.Range("L77").Value = 28933.5
$ws.Range("J77").Value = 9644.5
$ws.Range("N77").Value = -38293.5
$ws.Range("I81").Value = 1974.3572
$ws.Range("H81").Value = 5537.706
$ws.Range("K81").Value = 3948.7144
$ws.Range("M81").Value = -2887.7144
$ws.Range("I84").Value = 1974.3572
$ws.Range("H84").Value = 5537.706
$ws.Range("K84").Value = 19743.572
$ws.Range("M84").Value = -14439.572
$ws.Range("I107").Value = 828.63635
$ws.Range("H107").Value = 803.6799999999999
$ws.Range("K107").Value = 2485.90905
$ws.Range("M107").Value = -565.9090500000002
$ws.Range("L107").Value = 1862.0001
$ws.Range("J107").Value = 620.6667
$ws.Range("N107").Value = -5702.0001
$ws.Range("I132").Value = 1653.25
$ws.Range("H132").Value = 28991.7
$ws.Range("K132").Value = 4959.75
$ws.Range("M132").Value = -2429.75
$ws.Range("L132").Value = 168990.45
$ws.Range("J132").Value = 56330.15
$ws.Range("N132").Value = -174050.45
$ws.Range("I136").Value = 18544390
$ws.Range("H136").Value = 15264023
$ws.Range("K136").Value = 55633170
$ws.Range("M136").Value = -55630620
$ws.Range("L136").Value = 1507121.25
$ws.Range("J136").Value = 502373.75
$ws.Range("N136").Value = -1512221.25
